# Apply updated voltage-magnitude results for the 380 kV case (Case_4_28).
# Source data changed the slack/reference voltage set-points from 1.05 pu to 1.02 pu
# and the resulting bus voltage magnitudes were recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.10696588652695
$ws.Range("D2").Value = 1.103790675055804
$ws.Range("E2").Value = 1.118845340676255
$ws.Range("F2").Value = 1.121568073169455
$ws.Range("I2").Value = 1.081720913106794
$ws.Range("J2").Value = 1.111717419885484
$ws.Range("K2").Value = 1.10640445931726
$ws.Range("L2").Value = 1.121422090403256
$ws.Range("M2").Value = 1.124138239079049
$ws.Range("N2").Value = 1.113296186137492

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.108478968642181
$ws.Range("D3").Value = 1.10503660068001
$ws.Range("E3").Value = 1.120300600433435
$ws.Range("F3").Value = 1.122984388689149
$ws.Range("I3").Value = 1.082343583847421
$ws.Range("J3").Value = 1.112899011937741
$ws.Range("K3").Value = 1.107471934304521
$ws.Range("L3").Value = 1.122701109465684
$ws.Range("M3").Value = 1.12537887942661
$ws.Range("N3").Value = 1.114479456185993

# Row 4
$ws.Range("B4").Value = 1.019999999999999
$ws.Range("C4").Value = 1.109455724566376
$ws.Range("D4").Value = 1.105840245082653
$ws.Range("E4").Value = 1.121240186983429
$ws.Range("F4").Value = 1.12389869487757
$ws.Range("I4").Value = 1.082743845440841
$ws.Range("J4").Value = 1.113660765335917
$ws.Range("K4").Value = 1.108159539556617
$ws.Range("L4").Value = 1.123526121475147
$ws.Range("M4").Value = 1.126178969145299
$ws.Range("N4").Value = 1.115242291361369

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.109865809461273
$ws.Range("D5").Value = 1.106177493925995
$ws.Range("E5").Value = 1.121634703797054
$ws.Range("F5").Value = 1.124282564266744
$ws.Range("I5").Value = 1.08291148620298
$ws.Range("J5").Value = 1.113980340265914
$ws.Range("K5").Value = 1.108447868790934
$ws.Range("L5").Value = 1.12387234255061
$ws.Range("M5").Value = 1.126514691488388
$ws.Range("N5").Value = 1.115562320124423

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.109934632937402
$ws.Range("D6").Value = 1.106234084342918
$ws.Range("E6").Value = 1.121700916695494
$ws.Range("F6").Value = 1.124346988262704
$ws.Range("I6").Value = 1.082939597063999
$ws.Range("J6").Value = 1.114033959464463
$ws.Range("K6").Value = 1.108496237360358
$ws.Range("L6").Value = 1.123930438825946
$ws.Range("M6").Value = 1.126571023720533
$ws.Range("N6").Value = 1.115616015468384

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.109461206263725
$ws.Range("D7").Value = 1.105844753778848
$ws.Range("E7").Value = 1.121245460432781
$ws.Range("F7").Value = 1.123903826132632
$ws.Range("I7").Value = 1.082746087931434
$ws.Range("J7").Value = 1.11366503811972
$ws.Range("K7").Value = 1.108163395124406
$ws.Range("L7").Value = 1.123530750100443
$ws.Range("M7").Value = 1.126183457570319
$ws.Range("N7").Value = 1.115246570213015

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.107477723668662
$ws.Range("D8").Value = 1.104212274786645
$ws.Range("E8").Value = 1.119337585116112
$ws.Range("F8").Value = 1.122047172833906
$ws.Range("I8").Value = 1.08193189874433
$ws.Range("J8").Value = 1.112117332215126
$ws.Range("K8").Value = 1.10676586854504
$ws.Range("L8").Value = 1.121854883944944
$ws.Range("M8").Value = 1.124558080618436
$ws.Range("N8").Value = 1.113696666388512

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.103964455667371
$ws.Range("D9").Value = 1.101315731857145
$ws.Range("E9").Value = 1.115959466189741
$ws.Range("F9").Value = 1.118758703880666
$ws.Range("I9").Value = 1.080476699661862
$ws.Range("J9").Value = 1.109368157036057
$ws.Range("K9").Value = 1.104278995774293
$ws.Range("L9").Value = 1.118881520304149
$ws.Range("M9").Value = 1.121673013617593
$ws.Range("N9").Value = 1.110943587065354

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.101609489364928
$ws.Range("D10").Value = 1.099370819742478
$ws.Range("E10").Value = 1.113695943744549
$ws.Range("F10").Value = 1.116554551161839
$ws.Range("I10").Value = 1.079492488198634
$ws.Range("J10").Value = 1.107520136796509
$ws.Range("K10").Value = 1.102604298979394
$ws.Range("L10").Value = 1.116885127805522
$ws.Range("M10").Value = 1.119735044404867
$ws.Range("N10").Value = 1.109092942425095

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.100586594177711
$ws.Range("D11").Value = 1.098525246274916
$ws.Range("E11").Value = 1.112712975614104
$ws.Range("F11").Value = 1.115597199533686
$ws.Range("I11").Value = 1.079062907157615
$ws.Range("J11").Value = 1.106716197717674
$ws.Range("K11").Value = 1.101875049734831
$ws.Range("L11").Value = 1.116017195135485
$ws.Range("M11").Value = 1.118892307377779
$ws.Range("N11").Value = 1.108287861660556

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.10020615592443
$ws.Range("D12").Value = 1.098210640012821
$ws.Range("E12").Value = 1.112347418844308
$ws.Range("F12").Value = 1.115241144453077
$ws.Range("I12").Value = 1.078902823292175
$ws.Range("J12").Value = 1.106417008192202
$ws.Range("K12").Value = 1.101603549653959
$ws.Range("L12").Value = 1.115694273097265
$ws.Range("M12").Value = 1.118578728996673
$ws.Range("N12").Value = 1.10798824725164

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.100287783559762
$ws.Range("D13").Value = 1.098278147931714
$ws.Range("E13").Value = 1.112425852035663
$ws.Range("F13").Value = 1.115317540103134
$ws.Range("I13").Value = 1.078937185375695
$ws.Range("J13").Value = 1.106481211346486
$ws.Range("K13").Value = 1.101661815750409
$ws.Range("L13").Value = 1.115763565262902
$ws.Range("M13").Value = 1.118646017597391
$ws.Range("N13").Value = 1.108052541581767

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.100555157085982
$ws.Range("D14").Value = 1.09849925154138
$ws.Range("E14").Value = 1.112682767588227
$ws.Range("F14").Value = 1.115567777192727
$ws.Range("I14").Value = 1.07904968518354
$ws.Range("J14").Value = 1.106691478318141
$ws.Range("K14").Value = 1.101852620245945
$ws.Range("L14").Value = 1.115990513235178
$ws.Range("M14").Value = 1.118866398134487
$ws.Range("N14").Value = 1.10826310715664

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.10071982953017
$ws.Range("D15").Value = 1.098635411225871
$ws.Range("E15").Value = 1.112841003291653
$ws.Range("F15").Value = 1.115721896303843
$ws.Range("I15").Value = 1.079118931166236
$ws.Range("J15").Value = 1.106820954815111
$ws.Range("K15").Value = 1.101970098168882
$ws.Range("L15").Value = 1.116130272412007
$ws.Range("M15").Value = 1.119002108899403
$ws.Range("N15").Value = 1.108392767525087

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.101677307562886
$ws.Range("D16").Value = 1.099426864915441
$ws.Range("E16").Value = 1.113761119182181
$ws.Range("F16").Value = 1.116618024620884
$ws.Range("I16").Value = 1.079520925711646
$ws.Range("J16").Value = 1.107573412082079
$ws.Range("K16").Value = 1.102652609792965
$ws.Range("L16").Value = 1.116942655423243
$ws.Range("M16").Value = 1.119790897738358
$ws.Range("N16").Value = 1.109146293367681

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.102277048774934
$ws.Range("D17").Value = 1.099922402104445
$ws.Range("E17").Value = 1.11433751352169
$ws.Range("F17").Value = 1.117179348005921
$ws.Range("I17").Value = 1.079772168920801
$ws.Range("J17").Value = 1.108044402064794
$ws.Range("K17").Value = 1.103079628562561
$ws.Range("L17").Value = 1.117451302589558
$ws.Range("M17").Value = 1.120284717730454
$ws.Range("N17").Value = 1.109617952210195

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.102626561729409
$ws.Range("D18").Value = 1.100211111856211
$ws.Range("E18").Value = 1.114673440532749
$ws.Range("F18").Value = 1.117506475588988
$ws.Range("I18").Value = 1.079918386123905
$ws.Range("J18").Value = 1.108318762928385
$ws.Range("K18").Value = 1.103328306960869
$ws.Range("L18").Value = 1.11774765264777
$ws.Range("M18").Value = 1.120572409491115
$ws.Range("N18").Value = 1.109892702697683

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.102745685075506
$ws.Range("D19").Value = 1.10030949895989
$ws.Range("E19").Value = 1.114787936801705
$ws.Range("F19").Value = 1.117617969981082
$ws.Range("I19").Value = 1.079968186833498
$ws.Range("J19").Value = 1.108412252187434
$ws.Range("K19").Value = 1.103413033246492
$ws.Range("L19").Value = 1.117848643803744
$ws.Range("M19").Value = 1.12067044671327
$ws.Range("N19").Value = 1.109986324722202

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.102212733949659
$ws.Range("D20").Value = 1.099869269724046
$ws.Range("E20").Value = 1.114275700288685
$ws.Range("F20").Value = 1.117119152651637
$ws.Range("I20").Value = 1.079745246944738
$ws.Range("J20").Value = 1.107993906565222
$ws.Range("K20").Value = 1.103033854354315
$ws.Range("L20").Value = 1.117396764309949
$ws.Range("M20").Value = 1.120231771284484
$ws.Range("N20").Value = 1.109567385001223

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.100476435852724
$ws.Range("D21").Value = 1.098434156545121
$ws.Range("E21").Value = 1.112607124557015
$ws.Range("F21").Value = 1.115494101170865
$ws.Range("I21").Value = 1.079016571150323
$ws.Range("J21").Value = 1.106629575770173
$ws.Range("K21").Value = 1.10179645039041
$ws.Range("L21").Value = 1.115923697514708
$ws.Range("M21").Value = 1.118801516775925
$ws.Range("N21").Value = 1.108201116699954

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.099381918675029
$ws.Range("D22").Value = 1.097528815266636
$ws.Range("E22").Value = 1.111555481527888
$ws.Range("F22").Value = 1.114469746044366
$ws.Range("I22").Value = 1.078555421849973
$ws.Range("J22").Value = 1.105768459732584
$ws.Range("K22").Value = 1.101014827910708
$ws.Range("L22").Value = 1.114994431931211
$ws.Range("M22").Value = 1.117899081560147
$ws.Range("N22").Value = 1.10733877777882

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.099962415882054
$ws.Range("D23").Value = 1.098009044173164
$ws.Range("E23").Value = 1.112113222085683
$ws.Range("F23").Value = 1.115013028033319
$ws.Range("I23").Value = 1.078800172342692
$ws.Range("J23").Value = 1.106225270391864
$ws.Range("K23").Value = 1.10142952679817
$ws.Range("L23").Value = 1.115487349507625
$ws.Range("M23").Value = 1.11837778406924
$ws.Range("N23").Value = 1.107796237161633

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.10224179599623
$ws.Range("D24").Value = 1.099893278975456
$ws.Range("E24").Value = 1.114303631874878
$ws.Range("F24").Value = 1.117146353215253
$ws.Range("I24").Value = 1.079757412842596
$ws.Range("J24").Value = 1.1080167244192
$ws.Range("K24").Value = 1.103054538967483
$ws.Range("L24").Value = 1.117421408846839
$ws.Range("M24").Value = 1.12025569657073
$ws.Range("N24").Value = 1.10959023525917

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.104874928805633
$ws.Range("D25").Value = 1.102066967777391
$ws.Range("E25").Value = 1.116834767057989
$ws.Range("F25").Value = 1.11961089900434
$ws.Range("I25").Value = 1.080855362815084
$ws.Range("J25").Value = 1.110081533414651
$ws.Range("K25").Value = 1.104924834721541
$ws.Range("L25").Value = 1.119652661744362
$ws.Range("M25").Value = 1.122421406791084
$ws.Range("N25").Value = 1.111657976520231
